$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 and Row 4 swap identity: the "Gold Feb 26 / GC=F" row and the
# "Newmont Corporation / NEM" row trade places in the underlying shared
# string table, which (since the cell->string index mapping in the sheet
# itself did not change) means row 3 now displays Newmont/NEM and row 4
# now displays Gold Feb 26/GC=F.
$ws.Range("B3").Value = "Newmont Corporation"
$ws.Range("C3").Value = "NEM"
$ws.Range("B4").Value = "Gold Feb 26"
$ws.Range("C4").Value = "GC=F"

# --- Row 2 (StreetTRACKS Gold Shares / GLD) refreshed score columns
$ws.Range("K2").Value = 67.59999999999999
$ws.Range("N2").Value = 54.77309453746771

# --- Row 3 (now Newmont Corporation / NEM) refreshed data
$ws.Range("D3").Value = 90.72
$ws.Range("E3").Value = 52.7
$ws.Range("F3").Value = 0.22
$ws.Range("H3").Value = 80
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 86
$ws.Range("K3").Value = 66.40000000000001
$ws.Range("N3").Value = 54.77309453746771

# --- Row 4 (now Gold Feb 26 / GC=F) refreshed data
$ws.Range("D4").Value = 4204.1
$ws.Range("E4").Value = 52.1
$ws.Range("F4").Value = 0.93
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 76
$ws.Range("J4").Value = 83
$ws.Range("K4").Value = 64.8
$ws.Range("N4").Value = 54.77309453746771
